$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26..168 down to 27..169
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new record
$ws.Cells.Item(26, 1).Value = 8
$ws.Cells.Item(26, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44819
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = 100112044
$ws.Cells.Item(26, 7).Value = "Perejil"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 2000
$ws.Cells.Item(26, 11).Value = 2000
$ws.Cells.Item(26, 12).Value = 2500
$ws.Cells.Item(26, 13).Value = 2250
$ws.Cells.Item(26, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(26, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(26, 16).Value = 1500
$ws.Cells.Item(26, 17).Value = 1.5
$ws.Cells.Item(26, 18).Value = "Hortaliza"
